$wb = $excel.ActiveWorkbook

# Hunk @1465  Sheet=ALC  Row=17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2566.7778
$ws.Range("J17").Value = 2566.7778
$ws.Range("L17").Value = 7700.3334
$ws.Range("N17").Value = -8036.3334

# Hunk @4603  Sheet=ALC  Row=80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1472.375
$ws.Range("I80").Value = 5000
$ws.Range("K80").Value = 15000
$ws.Range("M80").Value = -14002

# Hunk @4750  Sheet=ALC  Row=83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1472.375
$ws.Range("I83").Value = 5000
$ws.Range("K83").Value = 45000
$ws.Range("M83").Value = -40008

# Hunk @4897  Sheet=ALC  Row=86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1835.6875
$ws.Range("I86").Value = 1852.8182
$ws.Range("J86").Value = 1798
$ws.Range("K86").Value = 1852.8182
$ws.Range("L86").Value = 1798
$ws.Range("M86").Value = -729.8181999999999
$ws.Range("N86").Value = -4044

# Hunk @5050  Sheet=ALC  Row=89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1835.6875
$ws.Range("I89").Value = 1852.8182
$ws.Range("J89").Value = 1798
$ws.Range("K89").Value = 9264.091
$ws.Range("L89").Value = 8990
$ws.Range("M89").Value = -3648.091
$ws.Range("N89").Value = -20222

# Hunk @5503  Sheet=ALC  Row=98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3108.5
$ws.Range("I98").Value = 3284.35
$ws.Range("K98").Value = 3284.35
$ws.Range("M98").Value = -1786.35

# Hunk @6703  Sheet=ALC  Row=122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 3108.5
$ws.Range("I122").Value = 3284.35
$ws.Range("K122").Value = 9853.05
$ws.Range("M122").Value = -7403.049999999999

# Hunk @7055  Sheet=ALC  Row=129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 791.7727
$ws.Range("J129").Value = 891.74286
$ws.Range("L129").Value = 2675.22858
$ws.Range("N129").Value = -12675.22858

# Hunk @8009  Sheet=ARM  Row=6
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 980
$ws.Range("I6").Value = 980
$ws.Range("K6").Value = 980
$ws.Range("M6").Value = -807

# Hunk @13648  Sheet=ARM  Row=122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1952.8422
$ws.Range("I122").Value = 1625.0667
$ws.Range("J122").Value = 3182
$ws.Range("K122").Value = 4875.2001
$ws.Range("L122").Value = 9546
$ws.Range("M122").Value = -2425.2001
$ws.Range("N122").Value = -14446

# Hunk @21115  Sheet=BSM  Row=134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1351.8
$ws.Range("I134").Value = 1104.8462
$ws.Range("K134").Value = 3314.5386
$ws.Range("M134").Value = -779.5385999999999

# Hunk @21806  Sheet=CRP  Row=6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1000
$ws.Range("I6").Value = 1000
$ws.Range("K6").Value = 1000
$ws.Range("M6").Value = -887

# Hunk @23019  Sheet=CRP  Row=31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1480.8572
$ws.Range("I31").Value = 1480.8572
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1480.8572
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1185.8572
$ws.Range("N31").ClearContents()

# Hunk @23169  Sheet=CRP  Row=34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1480.8572
$ws.Range("I34").Value = 1480.8572
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1480.8572
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1278.8572
$ws.Range("N34").ClearContents()

# Hunk @26324  Sheet=CRP  Row=99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1631.1111
$ws.Range("I99").Value = 1610
$ws.Range("K99").Value = 1610
$ws.Range("M99").Value = -112

# Hunk @27448  Sheet=CRP  Row=122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1156.8572
$ws.Range("I122").Value = 1183
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 3549
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -1099
$ws.Range("N122").Value = -7900

# Hunk @27638  Sheet=CRP  Row=126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1631.1111
$ws.Range("I126").Value = 1610
$ws.Range("K126").Value = 4830
$ws.Range("M126").Value = -2360

# Hunk @28635  Sheet=CUL  Row=4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3600142.5
$ws.Range("I4").Value = 2000185.2
$ws.Range("J4").Value = 8933333
$ws.Range("K4").Value = 6000555.6
$ws.Range("L4").Value = 26799999
$ws.Range("M4").Value = -6000443.6
$ws.Range("N4").Value = -26800223

# Hunk @30704  Sheet=CUL  Row=45
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 707.6667
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 707.6667
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 2123.0001
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -3187.0001

# Hunk @35104  Sheet=CUL  Row=131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J131").Value = 1665.7709
$ws.Range("L131").Value = 4997.3127
$ws.Range("N131").Value = -15077.3127

# Hunk @35618  Sheet=CUL  Row=141
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 58826130
$ws.Range("I141").Value = 76924370
$ws.Range("K141").Value = 230773110
$ws.Range("M141").Value = -230767930

# Hunk @41008  Sheet=GSM  Row=110
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 51000
$ws.Range("J110").Value = 51000
$ws.Range("L110").Value = 51000
$ws.Range("N110").Value = -59180

# Hunk @41590  Sheet=GSM  Row=122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2787.9092
$ws.Range("I122").Value = 2650.7856
$ws.Range("J122").Value = 3555.8
$ws.Range("K122").Value = 7952.3568
$ws.Range("L122").Value = 10667.4
$ws.Range("M122").Value = -5502.3568
$ws.Range("N122").Value = -15567.4

# Hunk @42927  Sheet=LTW  Row=7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2579.818
$ws.Range("I7").Value = 2263.111
$ws.Range("J7").Value = 4005
$ws.Range("K7").Value = 2263.111
$ws.Range("L7").Value = 4005
$ws.Range("M7").Value = -2151.111
$ws.Range("N7").Value = -4229

# Hunk @44523  Sheet=LTW  Row=40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5878.8184
$ws.Range("I40").Value = 2481.4
$ws.Range("J40").Value = 8710
$ws.Range("K40").Value = 2481.4
$ws.Range("L40").Value = 8710
$ws.Range("M40").Value = -2345.4
$ws.Range("N40").Value = -8982

# Hunk @45540  Sheet=LTW  Row=61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1750
$ws.Range("I61").Value = 1500
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1500
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1298
$ws.Range("N61").Value = -2404

# Hunk @48067  Sheet=LTW  Row=113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1750
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -6340

# Hunk @48698  Sheet=LTW  Row=126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2579.818
$ws.Range("I126").Value = 2263.111
$ws.Range("J126").Value = 4005
$ws.Range("K126").Value = 6789.333
$ws.Range("L126").Value = 12015
$ws.Range("M126").Value = -4319.333
$ws.Range("N126").Value = -16955

# Hunk @54913  Sheet=WVR  Row=113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 425.93332
$ws.Range("I113").Value = 359.625
$ws.Range("K113").Value = 1078.875
$ws.Range("M113").Value = 1091.125

# Hunk @55351  Sheet=WVR  Row=122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 11906740
$ws.Range("I122").Value = 13159949
$ws.Range("J122").Value = 1252.5
$ws.Range("K122").Value = 39479847
$ws.Range("L122").Value = 3757.5
$ws.Range("M122").Value = -39477397
$ws.Range("N122").Value = -8657.5

# Hunk @55550  Sheet=WVR  Row=126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 37037816
$ws.Range("I126").Value = 50000436
$ws.Range("J126").Value = 1760.5714
$ws.Range("K126").Value = 150001308
$ws.Range("L126").Value = 5281.7142
$ws.Range("M126").Value = -149998838
$ws.Range("N126").Value = -10221.7142
